$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Log In"
$ws.Range("B4").Value = "css"
$ws.Range("C4").Value = "input[class='btn'][value='Log In']"

$ws.Range("E10").Select()
